# Update gh-pages to output generated at 456a3b4
# Applies small numeric updates to the F column ("想去人数") on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Changes on "展览" sheet (F column)
$exhibitChanges = @{
    "F2"  = 14
    "F4"  = 402
    "F5"  = 5008
    "F6"  = 5008
    "F7"  = 63
    "F12" = 4799
    "F14" = 40
    "F19" = 237
    "F20" = 3707
    "F24" = 3517
    "F26" = 152
    "F36" = 6180
    "F37" = 975
    "F42" = 1278
    "F44" = 606
    "F45" = 23
    "F46" = 2141
    "F50" = 891
}

foreach ($cell in $exhibitChanges.Keys) {
    $wsExhibit.Range($cell).Value = $exhibitChanges[$cell]
}

# Changes on "全部类型" sheet (F column)
$allChanges = @{
    "F5"  = 402
    "F6"  = 5008
    "F7"  = 5008
    "F8"  = 63
    "F15" = 4799
    "F17" = 40
    "F22" = 237
    "F23" = 3707
    "F24" = 3517
    "F26" = 152
    "F35" = 6180
    "F36" = 975
    "F42" = 1278
    "F44" = 606
    "F45" = 2141
    "F49" = 891
}

foreach ($cell in $allChanges.Keys) {
    $wsAll.Range($cell).Value = $allChanges[$cell]
}
